$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New OneDrive links added for the Moffet et al. CosMx / Xenium rows ---
$ws.Range("E16").Value = "https://x2ocbw.bl.files.1drv.com/y4mR1ZXpUQbQIM5hSm5DuE0RmG5mync6IyU8Iodp_RFed9yo_ueW7pSlAfrR4DBt5AqodE6EOPlzaWCL-nkCGI2G5qaEs8mdUaKpQ_Ls0FUe20j7Bo4Jh_U10utmDdvOQDEuDWHXN22SamMQkeiJjgDZPk8_WQlaQ-z_7Him8CH_If5rFUjwh0N5Hu5escMoloId2NFHR3rVvf41f8Wtw6Z1Q"
$ws.Range("E17").Value = "https://9gwmhw.bl.files.1drv.com/y4m21lGsG84x5YvNlK65A5B6LkpxnX8NdK9AhhExllpnsfdRDHJla33U47c6opHZi6K0YoTJY54yriIsvhFU2t-D0zAZpjWmGCdoEMyjq2_F1JYsU6oqow5J5Dn0d_TDp9Cc2_shhh8OKqp-VvHSf9ihb4JNcMx9bE6bUxsUM7G6UEZqV9m9pMyTOQC9TN5JWHhKjMuPUsAo1TnBbxgJU3NYA"

# --- Row 99: author label corrected from "Vizgen" to "Wu et al." (matches rows 100-106) ---
$ws.Range("A99").Value = "Wu et al."

# --- Row 107 previously only had "Wu et al." in column A; replace it with a full
#     10X Genomics Visium dataset entry, then append six more similar rows (108-113) ---
$ws.Range("A107").Value = "10X"
$ws.Range("B107").Value = 2024
$ws.Range("D107").Value = "https://www.10xgenomics.com/datasets/adult-mouse-olfactory-bulb-1-standard-1"
$ws.Range("E107").Value = "https://kyhnwq.bl.files.1drv.com/y4mVei1sn6_3jISVzN2lgGJdq8JGByWfetd7S83TsBBEiBXjNuslaS55xKI4w3fIJAePL_putj1zH_Hl5AuaNn1epBqNhn6tsApjM_bMzn268fwpRDS4LNnzcbFUBN4SWYqIG4ue4maGC0bAbAsANZaXAXxy5tTdIjEn1VpboBvN6Y8Ds03djdxlEKigWwo45IP7soSzNSJ18NJaCws7nENZw"
$ws.Range("F107").Value = "Visium_olf_bulb"
$ws.Range("G107").Value = "Visium"
$ws.Range("H107").Value = "Olfactory bulb"
$ws.Range("I107").Value = "Mouse"
$ws.Range("M107").Value = 1

$ws.Range("A108").Value = "10X"
$ws.Range("B108").Value = 2024
$ws.Range("D108").Value = "https://www.10xgenomics.com/datasets/adult-mouse-kidney-ffpe-1-standard-1-3-0"
$ws.Range("E108").Value = "https://cg1oba.bl.files.1drv.com/y4mx4UjKNpxmTky9p84tsyymM9CIpQ_JqEu5PJOlc9G7iix7xuYQmmQPw_sVLgMLsSI3FjDNiQE8XIe_NW72eThWHxkdY-6PZBhiT6Z9kCRzhUfYQqLvB_tqMABSD02FqjX0bSvQ60jKX9hzZhVlD4EjM2z1bWWuhmiTPIpOFPQzECsRNH_NyuwM5JtV4BmRaNNdAhumA3dDJr0V7XtzRmj4w"
$ws.Range("F108").Value = "Visium_kidney_mouse"
$ws.Range("G108").Value = "Visium"
$ws.Range("H108").Value = "Kidney"
$ws.Range("I108").Value = "Mouse"
$ws.Range("M108").Value = 1

$ws.Range("A109").Value = "10X"
$ws.Range("B109").Value = 2024
$ws.Range("D109").Value = "https://www.10xgenomics.com/datasets/human-ovarian-cancer-1-standard"
$ws.Range("E109").Value = "https://mz23kq.bl.files.1drv.com/y4mG0545nd3mPv2J2pYWehC_UoqaWMTvR6PqZnpHLhTlPZ9PaM0rh2GdYrdAygLYVHLnkcW13KKgTl6L6QgzNCJvTOT16oSmRf49h-2y4Smq3hN3bSmG0pIb9rDL9FmOfV_S_-88gtQ3Xe6LiHF-dXb2Y5rXxDmqVZH5NuMKgSkE6MQlbiHPnJsaDlnK8oI0nw_3QYjW59kgSjKCDQVl4kSag"
$ws.Range("F109").Value = "Visium_ovarian_cancer"
$ws.Range("G109").Value = "Visium"
$ws.Range("H109").Value = "Ovarian cancer"
$ws.Range("I109").Value = "Human"
$ws.Range("M109").Value = 1

$ws.Range("A110").Value = "10X"
$ws.Range("B110").Value = 2024
$ws.Range("D110").Value = "https://www.10xgenomics.com/datasets/human-breast-cancer-visium-fresh-frozen-whole-transcriptome-1-standard"
$ws.Range("E110").Value = "https://mmgp7w.bl.files.1drv.com/y4m6rrHukO15I4Z0ayjDqsKOjLK67QXYeORSWY_P7N9OW9BFqj1cIchXQ9XYsEmuQ6a4s_e-7EzEq5k5ue2SFV1id3QXPKuF58QgIy6SPOCHP7hgpv2Kq3xQZFmqhlAMrynE1NY0UUQZyF288igTxIJbjDxrJgPShPzhRKhRp1S5YorjDnKvcykv02ICqqPJZ-8hc_2nbEPrAfHvjcM_mICKA"
$ws.Range("F110").Value = "Visium_breast_cancer"
$ws.Range("G110").Value = "Visium"
$ws.Range("H110").Value = "Breast cancer"
$ws.Range("I110").Value = "Human"
$ws.Range("M110").Value = 1

$ws.Range("A111").Value = "10X"
$ws.Range("B111").Value = 2024
$ws.Range("D111").Value = "https://www.10xgenomics.com/datasets/adult-mouse-brain-coronal-section-fresh-frozen-1-standard"
$ws.Range("E111").Value = "https://cpl1pg.bl.files.1drv.com/y4mDy2m7jRRx9lJZmCLQOl0y1RyG5s85iu_dOhJnU6ICMLOyZS5hhVe8pIM9pjb_zq_ZT05dIbEDAAoU418Pz1NgxgGXyHHGK6js-ho2whljFAEcF8OkBm-G1asnkXr721e4wqb60H6Ly9jlwD0hZC6vzp5VRSQw62ABh_8tdZHIAvn6yWU9YRvdiWlHSo3uvCcukTZoeAie8bq_Kfc9XqUsQ"
$ws.Range("F111").Value = "Visium_brain"
$ws.Range("G111").Value = "Visium"
$ws.Range("H111").Value = "Brain"
$ws.Range("I111").Value = "Mouse"
$ws.Range("M111").Value = 1

$ws.Range("A112").Value = "10X"
$ws.Range("B112").Value = 2024
$ws.Range("D112").Value = "https://www.10xgenomics.com/datasets/human-kidney-11-mm-capture-area-ffpe-2-standard"
$ws.Range("E112").Value = "https://ytqaqq.bl.files.1drv.com/y4mUcOiBHdjDdYz1DvTXv4Z1gNumtEZy-i4SCfaZrOuhbZJ0EVeuXx4dNrMbqwDLVYWNa6lhuWkzN1A7-Q52KpV84jbUaXVlCSv4gxB7P46RCT6t5BzJhAIqKTr7Mfy79FiT-khixVHPHpdhG0xtfgIqE3VdcPq698O8DycPOJj4AOyM1fJd4RjtYzXlDIUhowDiJcWt2T9QGorV2GtPCzPgQ"
$ws.Range("F112").Value = "Visium_kidney_human"
$ws.Range("G112").Value = "Visium"
$ws.Range("H112").Value = "Kidney"
$ws.Range("I112").Value = "Human"
$ws.Range("M112").Value = 1

$ws.Range("A113").Value = "10X"
$ws.Range("B113").Value = 2024
$ws.Range("D113").Value = "https://www.10xgenomics.com/datasets/visium-cytassist-mouse-embryo-11-mm-capture-area-ffpe-2-standard"
$ws.Range("E113").Value = "https://qvnota.bl.files.1drv.com/y4mkTerst2qgO8awvj5KruxDYL_4YI2UKXKRZ1JzDY_-9q38xh-j_oBLaNqiKvXL5X8SWV02OuE5WLaIXMQl1A0T_T3vtK25ZMxMpKKRjrVg_fReIIgKznneVW9rrOFciAInR37ESJ-RYeT1Ay5nRJWpYCpp_OCtOgz1CAl5Usn_OpDTd-_N4ubKID2ncO985fBWxn3hHoFyhIaMX0EoAgwaw"
$ws.Range("F113").Value = "Visium_embryo"
$ws.Range("G113").Value = "Visium"
$ws.Range("H113").Value = "Embryo"
$ws.Range("I113").Value = "Mouse"
$ws.Range("M113").Value = 1

# --- Window / selection state: scroll so row ~80 is at top, select D103 ---
$ws.Activate() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 80
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("D103").Select() | Out-Null

Write-Output "done"
